$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Save" column (H) mirroring the existing header style used by
# the other stat columns (B1:G1), then populate the new header + data cell.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

$ws.Range("H2").Value = 0
